# Apply the edits described by the commit:
# "Update data in excel spreadsheet to illustrate e2e pipeline"
#
# Summary of changes:
#  - GS sheet: update OM_Name/OM_Key/OM_Grade/OM_Team_Key values for rows 2 & 3
#  - GS sheet: widen column E a touch (cosmetic) and move selection
#  - CMS sheet: becomes the active/selected tab, with a new selection
#  - Workbook: CMS tab becomes the active tab (GS no longer active)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the data on the "GS" sheet (rows 2 and 3)
# ---------------------------------------------------------------------
$wsGS = $wb.Sheets.Item("GS")

# Row 2: James Smith / 1001 / PO / WMT  ->  Tom Swann / 1004 / C / WMT (ND02)
$wsGS.Range("E2").Value = "Tom Swann"
$wsGS.Range("F2").Value = 1004
$wsGS.Range("G2").Value = "C"
$wsGS.Range("H2").Value = "WMT (ND02)"

# Row 3: John Smith / 1002 / PO / WMT  ->  Andy Wright / 1005 / Z / WMT (ND02)
$wsGS.Range("E3").Value = "Andy Wright"
$wsGS.Range("F3").Value = 1005
$wsGS.Range("G3").Value = "Z"
$wsGS.Range("H3").Value = "WMT (ND02)"

# Slightly widen column E on the GS sheet
$wsGS.Columns.Item(5).ColumnWidth = 14.3

# ---------------------------------------------------------------------
# 2. Switch the active sheet/selection from "GS" to "CMS"
# ---------------------------------------------------------------------
# Leave a new resting selection on the GS sheet first (selecting on it
# later would re-activate it, so do this before switching tabs)
$wsGS.Range("N21").Select() | Out-Null

$wsCMS = $wb.Sheets.Item("CMS")
$wsCMS.Activate() | Out-Null
$wsCMS.Range("K13").Select() | Out-Null
